$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "60.729.61"
$ws.Cells.Item(2, 5).Value = "  -1.22%  "
$ws.Cells.Item(3, 4).Value = "2.906.30"
$ws.Cells.Item(3, 5).Value = "  -1.77%  "
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
$ws.Cells.Item(5, 4).Value = "'529.04"
$ws.Cells.Item(5, 5).Value = "  -2.48%  "
$ws.Cells.Item(6, 4).Value = "'144.05"
$ws.Cells.Item(6, 5).Value = "  -5.16%  "
$ws.Cells.Item(7, 4).Value = "'0.998"
$ws.Cells.Item(7, 5).Value = "  -0.11%  "
$ws.Cells.Item(8, 4).Value = "'0.556"
$ws.Cells.Item(8, 5).Value = "  -2.39%  "
$ws.Cells.Item(9, 4).Value = "2.916.58"
$ws.Cells.Item(9, 5).Value = "  -1.70%  "
$ws.Cells.Item(10, 5).Value = "  -4.13%  "
$ws.Cells.Item(11, 4).Value = "'6.06"
$ws.Cells.Item(11, 5).Value = "  -1.05%  "
$ws.Cells.Item(12, 4).Value = "'0.362"
$ws.Cells.Item(12, 5).Value = "  -1.78%  "
$ws.Cells.Item(13, 4).Value = "3.416.19"
$ws.Cells.Item(13, 5).Value = "  -1.71%  "
$ws.Cells.Item(14, 5).Value = "  +1.75%  "
$ws.Cells.Item(15, 4).Value = "60.668.25"
$ws.Cells.Item(15, 5).Value = "  -1.45%  "
$ws.Cells.Item(16, 4).Value = "'22.88"
$ws.Cells.Item(16, 5).Value = "  -3.45%  "
$ws.Cells.Item(17, 4).Value = "2.911.43"
$ws.Cells.Item(17, 5).Value = "  -1.65%  "
$ws.Cells.Item(18, 4).Value = "'0.0000142"
$ws.Cells.Item(18, 5).Value = "  -3.49%  "
$ws.Cells.Item(19, 4).Value = "'5.05"
$ws.Cells.Item(19, 5).Value = "  -2.30%  "
$ws.Cells.Item(20, 4).Value = "'11.76"
$ws.Cells.Item(20, 5).Value = "  -1.89%  "
$ws.Cells.Item(21, 4).Value = "'362.92"
$ws.Cells.Item(21, 5).Value = "  -4.99%  "
$ws.Cells.Item(22, 4).Value = "'6.64"
$ws.Cells.Item(22, 5).Value = "  -0.50%  "
$ws.Cells.Item(23, 5).Value = "  +0.00%  "
$ws.Cells.Item(24, 4).Value = "'5.68"
$ws.Cells.Item(24, 5).Value = "  +0.28%  "
$ws.Cells.Item(25, 4).Value = "'64.76"
$ws.Cells.Item(25, 5).Value = "  -0.82%  "
$ws.Cells.Item(26, 4).Value = "'0.456"
$ws.Cells.Item(26, 5).Value = "  -2.72%  "
$ws.Cells.Item(27, 4).Value = "'0.180"
$ws.Cells.Item(27, 5).Value = "  -3.42%  "
$ws.Cells.Item(28, 4).Value = "'0.997"
$ws.Cells.Item(28, 5).Value = "  -0.10%  "
$ws.Cells.Item(29, 4).Value = "'7.87"
$ws.Cells.Item(29, 5).Value = "  -5.68%  "
$ws.Cells.Item(30, 4).Value = "0.0₃0858"
$ws.Cells.Item(30, 5).Value = "  -8.10%  "
$ws.Cells.Item(31, 4).Value = "'1.00"
$ws.Cells.Item(31, 5).Value = "  +0.07%  "
$ws.Cells.Item(32, 4).Value = "'1.69"
$ws.Cells.Item(32, 5).Value = "  -2.17%  "
$ws.Cells.Item(33, 4).Value = "'19.81"
$ws.Cells.Item(33, 5).Value = "  -3.06%  "
$ws.Cells.Item(34, 4).Value = "'152.30"
$ws.Cells.Item(34, 5).Value = "  -4.49%  "
$ws.Cells.Item(35, 4).Value = "'4.38"
$ws.Cells.Item(35, 5).Value = "  -5.98%  "
$ws.Cells.Item(36, 4).Value = "'5.60"
$ws.Cells.Item(36, 5).Value = "  -6.01%  "
$ws.Cells.Item(37, 4).Value = "'1.01"
$ws.Cells.Item(37, 5).Value = "  -5.34%  "
$ws.Cells.Item(38, 5).Value = "  -4.81%  "
$ws.Cells.Item(39, 4).Value = "'37.84"
$ws.Cells.Item(39, 5).Value = "  +1.68%  "
$ws.Cells.Item(40, 4).Value = "'1.49"
$ws.Cells.Item(40, 5).Value = "  -4.02%  "
$ws.Cells.Item(41, 4).Value = "'3.73"
$ws.Cells.Item(41, 5).Value = "  -5.16%  "
$ws.Cells.Item(42, 4).Value = "2.297.22"
$ws.Cells.Item(42, 5).Value = "  -4.77%  "
$ws.Cells.Item(43, 4).Value = "'0.650"
$ws.Cells.Item(43, 5).Value = "  -1.87%  "
$ws.Cells.Item(44, 4).Value = "'0.0585"
$ws.Cells.Item(44, 5).Value = "  -1.70%  "
$ws.Cells.Item(45, 4).Value = "'20.50"
$ws.Cells.Item(45, 5).Value = "  -7.48%  "
$ws.Cells.Item(47, 4).Value = "'5.01"
$ws.Cells.Item(47, 5).Value = "  +0.73%  "
$ws.Cells.Item(48, 5).Value = "  -3.15%  "
$ws.Cells.Item(49, 4).Value = "'10.32"
$ws.Cells.Item(49, 5).Value = "  -1.30%  "
$ws.Cells.Item(50, 4).Value = "'0.0924"
$ws.Cells.Item(50, 5).Value = "  -3.48%  "
$ws.Cells.Item(51, 4).Value = "'251.77"
$ws.Cells.Item(51, 5).Value = "  -6.13%  "
